$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Columns B-E for rows 2-51 hold text-like values (including numeric-looking
# strings such as "217.44"); force text format before/after writing so COM
# does not coerce them into numeric cells, then strip the format again so
# the cells keep their original (default) style.
$dataRange = $ws.Range("B2:E51")
$dataRange.NumberFormat = "@"

$ws.Range("D2").Value = "26.107.77"
$ws.Range("E2").Value = "  -0.80%  "
$ws.Range("D3").Value = "1.649.90"
$ws.Range("E3").Value = "  -0.95%  "
$ws.Range("E4").Value = "  -0.49%  "
$ws.Range("D5").Value = "217.44"
$ws.Range("E5").Value = "  -0.81%  "
$ws.Range("D6").Value = "0.5195"
$ws.Range("E6").Value = "  -2.99%  "
$ws.Range("E7").Value = "  -0.44%  "
$ws.Range("D8").Value = "0.2617"
$ws.Range("E8").Value = "  -1.64%  "
$ws.Range("E9").Value = "  -2.22%  "
$ws.Range("D10").Value = "20.41"
$ws.Range("E10").Value = "  -1.61%  "
$ws.Range("D11").Value = "0.07793"
$ws.Range("E11").Value = "  -0.72%  "
$ws.Range("D12").Value = "1.680.84"
$ws.Range("E12").Value = "  +0.57%  "
$ws.Range("E13").Value = "  -2.29%  "
$ws.Range("D14").Value = "1.876.87"
$ws.Range("E14").Value = "  -0.87%  "
$ws.Range("D15").Value = "0.5544"
$ws.Range("E15").Value = "  +0.07%  "
$ws.Range("D16").Value = "0.0₅7976"
$ws.Range("E16").Value = "  -3.26%  "
$ws.Range("D17").Value = "64.69"
$ws.Range("E17").Value = "  -1.71%  "
$ws.Range("D18").Value = "26.099.71"
$ws.Range("E18").Value = "  -0.89%  "
$ws.Range("E19").Value = "  -0.44%  "
$ws.Range("D20").Value = "4.625"
$ws.Range("E20").Value = "  -1.72%  "
$ws.Range("D21").Value = "193.94"
$ws.Range("E21").Value = "  +0.08%  "
$ws.Range("D22").Value = "10.06"
$ws.Range("E22").Value = "  -2.17%  "
$ws.Range("D23").Value = "5.943"
$ws.Range("E23").Value = "  -1.62%  "
$ws.Range("E24").Value = "  -0.56%  "
$ws.Range("D25").Value = "146.76"
$ws.Range("E25").Value = "  +0.30%  "
$ws.Range("E26").Value = "  -2.33%  "
$ws.Range("D27").Value = "7.161"
$ws.Range("E27").Value = "  -0.55%  "
$ws.Range("D28").Value = "15.92"
$ws.Range("E28").Value = "  -1.33%  "
$ws.Range("E29").Value = "  -0.67%  "
$ws.Range("D30").Value = "0.05621"
$ws.Range("E30").Value = "  -3.70%  "
$ws.Range("D31").Value = "1.265"
$ws.Range("E31").Value = "  -1.39%  "
$ws.Range("D32").Value = "3.484"
$ws.Range("E32").Value = "  -3.74%  "
$ws.Range("D33").Value = "3.375"
$ws.Range("E33").Value = "  +2.80%  "
$ws.Range("D34").Value = "1.596"
$ws.Range("E34").Value = "  -1.40%  "
$ws.Range("D35").Value = "2.801"
$ws.Range("E35").Value = "  -0.86%  "
$ws.Range("D36").Value = "0.9469"
$ws.Range("E36").Value = "  -2.49%  "
$ws.Range("E37").Value = "  -0.73%  "
$ws.Range("D38").Value = "0.5645"
$ws.Range("E38").Value = "  -3.14%  "
$ws.Range("D39").Value = "5.955"
$ws.Range("E39").Value = "  +1.49%  "
$ws.Range("D40").Value = "0.01577"
$ws.Range("E40").Value = "  -1.83%  "
$ws.Range("D41").Value = "1.060.49"
$ws.Range("E41").Value = "  +0.72%  "
$ws.Range("D43").Value = "0.8378"
$ws.Range("E43").Value = "  -3.97%  "
$ws.Range("D44").Value = "102.97"
$ws.Range("E44").Value = "  -2.20%  "
$ws.Range("B45").Value = "RocketPoolETH"
$ws.Range("C45").Value = "https://coinranking.com/coin/QJZRUGyNI+rocketpooleth-reth"
$ws.Range("D45").Value = "1.787.59"
$ws.Range("E45").Value = "  -0.93%  "
$ws.Range("B46").Value = "Aave"
$ws.Range("C46").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Range("D46").Value = "57.04"
$ws.Range("E46").Value = "  -1.45%  "
$ws.Range("B47").Value = "BabyDogeCoin"
$ws.Range("C47").Value = "https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge"
$ws.Range("D47").Value = "0.0₈106"
$ws.Range("E47").Value = "  -0.74%  "
$ws.Range("B48").Value = "Frax"
$ws.Range("C48").Value = "https://coinranking.com/coin/KfWtaeV1W+frax-frax"
$ws.Range("D48").Value = "1.005"
$ws.Range("E48").Value = "  -0.94%  "
$ws.Range("B49").Value = "Cronos"
$ws.Range("C49").Value = "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
$ws.Range("D49").Value = "0.05314"
$ws.Range("E49").Value = "  +2.86%  "
$ws.Range("B50").Value = "Mantle"
$ws.Range("C50").Value = "https://coinranking.com/coin/BoI4ux0nd+mantle-mnt"
$ws.Range("D50").Value = "0.4333"
$ws.Range("E50").Value = "  -1.16%  "
$ws.Range("B51").Value = "EnergySwap"
$ws.Range("C51").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D51").Value = "7.950"
$ws.Range("E51").Value = "  -1.35%  "

$dataRange.ClearFormats()
